# Fix xlsx photometer data upload:
#  - the worksheet was mis-named "photometer_xlsv" (a stray extension-like
#    typo) when it should be "photometer_data", matching the workbook file.
#  - reset the saved cursor/selection back to the top-left cell (A1)
#    instead of the stray A4 selection left over from editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (typo fix).
$ws.Name = "photometer_data"

# Reset selection to A1.
$ws.Range("A1").Select()
